$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Make room for two new columns (tahun_akademik, prodi) between
#        "pekerjaan_ortu" (F) and the old "nama_prestasi" (G) by moving the
#        existing G:L block two columns to the right (I:N). Done column by
#        column, right-to-left, so source/destination never overlap.
$ws.Range("L1:L9").Cut($ws.Range("N1"))
$ws.Range("K1:K9").Cut($ws.Range("M1"))
$ws.Range("J1:J9").Cut($ws.Range("L1"))
$ws.Range("I1:I9").Cut($ws.Range("K1"))
$ws.Range("H1:H9").Cut($ws.Range("J1"))
$ws.Range("G1:G9").Cut($ws.Range("I1"))

# --- 2) Fix the header typo that rode along with the move (L1).
$ws.Range("L1").Value = "jenis_sertifikat"

# --- 3) Make sure the header row is fully bold/styled like the rest (B1
#        had been missed previously).
$ws.Range("B1").Font.Bold = $true

# --- 4) prodi column should store its values as text ("04", not 4).
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H7").NumberFormat = "@"

# --- 5) New column headers.
$ws.Range("G1").Value = "tahun_akademik"
$ws.Range("H1").Value = "prodi"
$ws.Range("H1").NumberFormat = "@"
$ws.Range("H1").Font.Bold = $true

# --- 6) New data for the three student blocks.
$ws.Range("G2").Value = 2019
$ws.Range("H2").Value = "04"

$ws.Range("G4").Value = 2019
$ws.Range("H4").Value = "04"

$ws.Range("G7").Value = 2019
$ws.Range("H7").Value = "04"

# --- 7) Renumber the nim (student id) column.
$ws.Range("A2").Value = 19090001
$ws.Range("A3").Value = 19090001
$ws.Range("A4").Value = 19090002
$ws.Range("A5").Value = 19090002
$ws.Range("A6").Value = 19090002
$ws.Range("A7").Value = 19090003
$ws.Range("A8").Value = 19090003
$ws.Range("A9").Value = 19090003

# --- 8) Column width tweaks that came with the rework.
$ws.Columns("K").ColumnWidth = 24.21875

# --- 9) Tidy up the empty-husk cells that the column shuffle (step 1)
#        leaves behind in rows that don't use every column.
$ws.Range("H3:L3").ClearContents()
$ws.Range("H5:J5").ClearContents()
$ws.Range("M5:N5").ClearContents()
$ws.Range("H6:J6").ClearContents()
$ws.Range("M6:N6").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("H9").ClearContents()

# --- 10) Selection, matching the saved state in the workbook.
$ws.Range("N7").Select()
